$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4424.5
$ws.Range("I2").Value = 3187.5
$ws.Range("K2").Value = 3187.5
$ws.Range("M2").Value = -3074.5

$ws.Range("H33").Value = 16862496
$ws.Range("I33").Value = 7268838
$ws.Range("K33").Value = 7268838
$ws.Range("M33").Value = -7268609

$ws.Range("H51").Value = 83503900
$ws.Range("J51").Value = 166670400
$ws.Range("L51").Value = 166670400
$ws.Range("N51").Value = -166671368

$ws.Range("H58").Value = 617.6
$ws.Range("I58").Value = 617.6
$ws.Range("K58").Value = 1852.8
$ws.Range("M58").Value = -1702.8

$ws.Range("H82").Value = 2566.7144
$ws.Range("I82").Value = 2566.7144
$ws.Range("K82").Value = 7700.1432
$ws.Range("M82").Value = -7294.1432

$ws.Range("H85").Value = 2566.7144
$ws.Range("I85").Value = 2566.7144
$ws.Range("K85").Value = 7700.1432
$ws.Range("M85").Value = -6296.1432

$ws.Range("H104").Value = 1333.3334
$ws.Range("I104").Value = 1000
$ws.Range("K104").Value = 3000
$ws.Range("M104").Value = -1253

$ws.Range("H135").Value = 645.9231
$ws.Range("I135").Value = 550.7778
$ws.Range("K135").Value = 4957.000199999999
$ws.Range("M135").Value = -2422.000199999999

$ws.Range("H137").Value = 5011234
$ws.Range("I137").Value = 8975.77
$ws.Range("K137").Value = 26927.31
$ws.Range("M137").Value = -24377.31

$ws.Range("H141").Value = 3819.7273
$ws.Range("I141").Value = 1335.2222
$ws.Range("J141").Value = 15000
$ws.Range("K141").Value = 4005.6666
$ws.Range("L141").Value = 45000
$ws.Range("M141").Value = 1174.3334
$ws.Range("N141").Value = -55360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2049.5
$ws.Range("I2").Value = 2000
$ws.Range("K2").Value = 2000
$ws.Range("M2").Value = -1887

$ws.Range("H45").Value = 60309
$ws.Range("I45").Value = 72810.42999999999
$ws.Range("J45").Value = 1969
$ws.Range("K45").Value = 72810.42999999999
$ws.Range("L45").Value = 1969
$ws.Range("M45").Value = -72433.42999999999
$ws.Range("N45").Value = -2723

$ws.Range("H61").Value = 1278536.9
$ws.Range("I61").Value = 39925.586
$ws.Range("K61").Value = 39925.586
$ws.Range("M61").Value = -39713.586

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H97").Value = 3620.4146
$ws.Range("I97").Value = 3679.182
$ws.Range("K97").Value = 3679.182
$ws.Range("M97").Value = -3183.182

$ws.Range("H116").Value = 2049.5
$ws.Range("I116").Value = 2000
$ws.Range("K116").Value = 2000
$ws.Range("M116").Value = 294

$ws.Range("H132").Value = 2571.75
$ws.Range("I132").Value = 2109.9546
$ws.Range("J132").Value = 3587.7
$ws.Range("K132").Value = 6329.8638
$ws.Range("L132").Value = 10763.1
$ws.Range("M132").Value = -3799.8638
$ws.Range("N132").Value = -15823.1

$ws.Range("H136").Value = 1278536.9
$ws.Range("I136").Value = 39925.586
$ws.Range("K136").Value = 119776.758
$ws.Range("M136").Value = -117226.758

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2049.5
$ws.Range("I3").Value = 2000
$ws.Range("K3").Value = 2000
$ws.Range("M3").Value = -1886

$ws.Range("H80").Value = 778.7273
$ws.Range("I80").Value = 1065.6666
$ws.Range("K80").Value = 1065.6666
$ws.Range("M80").Value = -67.66660000000002

$ws.Range("H83").Value = 778.7273
$ws.Range("I83").Value = 1065.6666
$ws.Range("K83").Value = 5328.333000000001
$ws.Range("M83").Value = -336.3330000000005

$ws.Range("H94").Value = 1288.2572
$ws.Range("I94").Value = 1384.7142
$ws.Range("J94").Value = 902.4286
$ws.Range("K94").Value = 1384.7142
$ws.Range("L94").Value = 902.4286
$ws.Range("M94").Value = -933.7141999999999
$ws.Range("N94").Value = -1804.4286

$ws.Range("H107").Value = 32552.375
$ws.Range("I107").Value = 38403.332
$ws.Range("K107").Value = 38403.332
$ws.Range("M107").Value = -36483.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4186.2256
$ws.Range("I31").Value = 3364.5625
$ws.Range("K31").Value = 3364.5625
$ws.Range("M31").Value = -3069.5625

$ws.Range("H34").Value = 4186.2256
$ws.Range("I34").Value = 3364.5625
$ws.Range("K34").Value = 3364.5625
$ws.Range("M34").Value = -3162.5625

$ws.Range("H58").Value = 1421.7307
$ws.Range("I58").Value = 1186
$ws.Range("K58").Value = 1186
$ws.Range("M58").Value = -983

$ws.Range("H62").Value = 4424.3076
$ws.Range("I62").Value = 3551.4
$ws.Range("J62").Value = 4969.875
$ws.Range("K62").Value = 3551.4
$ws.Range("L62").Value = 4969.875
$ws.Range("M62").Value = -2927.4
$ws.Range("N62").Value = -6217.875

$ws.Range("H65").Value = 4424.3076
$ws.Range("I65").Value = 3551.4
$ws.Range("J65").Value = 4969.875
$ws.Range("K65").Value = 17757
$ws.Range("L65").Value = 24849.375
$ws.Range("M65").Value = -14637
$ws.Range("N65").Value = -31089.375

$ws.Range("H136").Value = 1421.7307
$ws.Range("I136").Value = 1186
$ws.Range("K136").Value = 3558
$ws.Range("M136").Value = -1008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 756
$ws.Range("I8").Value = 756
$ws.Range("K8").Value = 2268
$ws.Range("M8").Value = -2129

$ws.Range("H50").Value = 9180.5
$ws.Range("J50").Value = 11560
$ws.Range("L50").Value = 34680
$ws.Range("N50").Value = -35642

$ws.Range("H53").Value = 9180.5
$ws.Range("J53").Value = 11560
$ws.Range("L53").Value = 34680
$ws.Range("N53").Value = -35642

$ws.Range("H134").Value = 2831.923
$ws.Range("I134").Value = 710.4545000000001
$ws.Range("K134").Value = 2131.3635
$ws.Range("M134").Value = 2938.6365

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5751.8
$ws.Range("I70").Value = 5586.8335
$ws.Range("K70").Value = 5586.8335
$ws.Range("M70").Value = -5316.8335

$ws.Range("H73").Value = 5751.8
$ws.Range("I73").Value = 5586.8335
$ws.Range("K73").Value = 5586.8335
$ws.Range("M73").Value = -4650.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6949006.5
$ws.Range("I40").Value = 9263342
$ws.Range("K40").Value = 9263342
$ws.Range("M40").Value = -9263206

$ws.Range("H122").Value = 2790.1
$ws.Range("I122").Value = 2363
$ws.Range("J122").Value = 3074.8333
$ws.Range("K122").Value = 7089
$ws.Range("L122").Value = 9224.499899999999
$ws.Range("M122").Value = -4639
$ws.Range("N122").Value = -14124.4999

$ws.Range("H132").Value = 5332.0835
$ws.Range("I132").Value = 4573.5713
$ws.Range("K132").Value = 13720.7139
$ws.Range("M132").Value = -11190.7139

$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 50000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 50000
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1101152.6
$ws.Range("I107").Value = 2160.6155
$ws.Range("J107").Value = 2200144.5
$ws.Range("K107").Value = 6481.8465
$ws.Range("L107").Value = 6600433.5
$ws.Range("M107").Value = -4561.8465
$ws.Range("N107").Value = -6604273.5

$ws.Range("H113").Value = 442.4762
$ws.Range("I113").Value = 381.42856
$ws.Range("K113").Value = 1144.28568
$ws.Range("M113").Value = 1025.71432

$ws.Range("H122").Value = 5195.909
$ws.Range("I122").Value = 5115.5
$ws.Range("K122").Value = 15346.5
$ws.Range("M122").Value = -12896.5

$ws.Range("H132").Value = 1905.7322
$ws.Range("I132").Value = 1767
$ws.Range("K132").Value = 5301
$ws.Range("M132").Value = -2771
